$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.881.77"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.627.84"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'211.37"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'0.522"
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'23.48"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.858.95"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "1.621.21"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'0.563"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'65.55"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "27.862.45"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'230.66"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'7.66"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "'154.65"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'15.55"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D34").Value = "1.399.99"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").Value = "'0.557"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'1.84"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "'65.92"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "1.769.45"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "'88.12"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.103"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  -0.48%  "
